$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # The source cells store these look-like-numbers/dates values as plain
    # text (inlineStr). Assigning the string directly would let Excel's
    # auto-conversion turn "2026/01/20" into a date serial or "8.88" into a
    # number, so force the Text number format for the assignment and then
    # drop back to the Normal style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Rows whose Date_1 (column A) advances from 2026/01/19 to 2026/01/20.
$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $dateRows) {
    Set-TextValue $ws.Range("A$r") "2026/01/20"
}

# EBITDA (column B) value refreshes alongside the date on a few rows.
Set-TextValue $ws.Range("B8")  "8.88"
Set-TextValue $ws.Range("B14") "3.20"
Set-TextValue $ws.Range("B56") "32.76"
